# Updates cryptocurrency price (D) and 1h volume-change (E) columns
# to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Some "Price" cells hold digit-only text (e.g. "308.24") that Excel would
# otherwise auto-convert to a number on assignment. Briefly mark just those
# cells as Text first so the values round-trip as the exact original strings
# (preserving trailing zeros etc.), then drop the temporary formatting again
# so the cell style stays exactly as it was.
$forceTextCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D14", "D15", "D17", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D35", "D37", "D38", "D39", "D40", "D41", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "308.24"
$ws.Range("D6").Value = "83.15"
$ws.Range("D7").Value = "0.525"
$ws.Range("D9").Value = "0.476"
$ws.Range("D10").Value = "0.0797"
$ws.Range("D11").Value = "29.49"
$ws.Range("D14").Value = "6.34"
$ws.Range("D15").Value = "14.58"
$ws.Range("D17").Value = "0.748"
$ws.Range("D20").Value = "6.00"
$ws.Range("D21").Value = "67.69"
$ws.Range("D22").Value = "10.40"
$ws.Range("D23").Value = "233.12"
$ws.Range("D24").Value = "2.51"
$ws.Range("D26").Value = "1.79"
$ws.Range("D27").Value = "23.23"
$ws.Range("D28").Value = "2.19"
$ws.Range("D29").Value = "9.16"
$ws.Range("D30").Value = "33.56"
$ws.Range("D31").Value = "152.33"
$ws.Range("D33").Value = "5.02"
$ws.Range("D35").Value = "0.0709"
$ws.Range("D37").Value = "2.73"
$ws.Range("D38").Value = "0.0967"
$ws.Range("D39").Value = "15.27"
$ws.Range("D40").Value = "1.68"
$ws.Range("D41").Value = "3.71"
$ws.Range("D44").Value = "0.0261"
$ws.Range("D45").Value = "17.14"
$ws.Range("D46").Value = "9.37"
$ws.Range("D47").Value = "2.64"
$ws.Range("D49").Value = "91.74"
$ws.Range("D50").Value = "69.38"
$ws.Range("D51").Value = "49.19"

foreach ($addr in $forceTextCells) {
    $ws.Range($addr).ClearFormats()
}

# Remaining Price/Volume cells are safe to set directly (either they
# contain a "%" sign, or the text has multiple "." separators so Excel
# keeps it as text automatically).
$ws.Range("D2").Value = "39.763.69"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "2.325.63"
$ws.Range("E3").Value = "  -4.14%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("E6").Value = "  -6.94%  "
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("E10").Value = "  -4.45%  "
$ws.Range("E11").Value = "  -7.97%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "2.693.70"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("E14").Value = "  -5.83%  "
$ws.Range("E15").Value = "  -6.77%  "
$ws.Range("D16").Value = "2.350.10"
$ws.Range("E16").Value = "  -4.02%  "
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "39.706.27"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "0.0₃0891"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("E20").Value = "  -4.80%  "
$ws.Range("E21").Value = "  -6.34%  "
$ws.Range("E22").Value = "  -5.65%  "
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("E24").Value = "  -6.66%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("E27").Value = "  -3.77%  "
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").Value = "  -4.84%  "
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -5.03%  "
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").Value = "  -7.26%  "
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("E39").Value = "  -9.95%  "
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("D42").Value = "1.965.59"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("E44").Value = "  -5.49%  "
$ws.Range("E45").Value = "  -8.05%  "
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("E47").Value = "  -9.15%  "
$ws.Range("D48").Value = "2.556.10"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("E51").Value = "  -5.31%  "
